$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 107, shifting existing rows 107-120 down to 108-121.
$ws.Rows.Item(107).Insert()

# New row 107 keeps the same market/category/unit info as the row that was
# previously at 107 (now at 108), but records a new weekly observation.
$ws.Cells.Item(107, 1).Value = 10
$ws.Cells.Item(107, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(107, 3).Value = "La Araucanía"
$ws.Cells.Item(107, 4).Value = 44449
$ws.Cells.Item(107, 4).NumberFormat = $ws.Cells.Item(108, 4).NumberFormat
$ws.Cells.Item(107, 5).Value = 9
$ws.Cells.Item(107, 6).Value = 100112005
$ws.Cells.Item(107, 7).Value = "Puerro"
$ws.Cells.Item(107, 8).Value = "Azul de Maquehue"
$ws.Cells.Item(107, 9).Value = "Primera"
$ws.Cells.Item(107, 10).Value = 65
$ws.Cells.Item(107, 11).Value = 8000
$ws.Cells.Item(107, 12).Value = 8000
$ws.Cells.Item(107, 13).Value = 8000
$ws.Cells.Item(107, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(107, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(107, 16).Value = 667
$ws.Cells.Item(107, 17).Value = 12
$ws.Cells.Item(107, 18).Value = "Hortaliza"
